$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution time for the "html css" row (row 4) from "0h 5m" to "0h 22m"
$ws.Range("C4").Value = "0h 22m"

# Update the active selection to C5 (matches the saved selection state in the file)
$ws.Range("C5").Select()
